# Apply the "add S/T review-opinion columns" edit to the '2. Details' sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2. Details")

# All data rows (3..387) that have a test-name value in column R and need
# column S ("review result 1") / column T ("review result 2") populated.
$rows = @(3,4,5,6,7,8,9,10,11,12,14,17,18,19,20,21,22,23,24,25,26,27,28,29,30,31,32,33,
50,51,52,56,57,
69,70,71,72,73,74,75,76,77,78,79,80,81,82,83,
259,260,261,262,263,264,265,266,267,268,269,270,271,272,273,274,275,276,277,278,279,280,281,282,283,
303,304,305,306,307,308,309,310,311,312,313,314,315,316,317,318,319,320,321,322,323,324,325,326,327,328,329,330,
331,332,333,334,335,336,337,338,339,340,341,342,343,344,345,346,347,348,349,350,
351,352,353,354,355,356,357,358,359,360,361,362,363,364,365,366,367,368,369,370,
371,372,373,374,375,376,377,378,379,380,381,382,383,384,385,386,387)

# Rows where column S gets a specific test-function name instead of the
# default "x" marker.
$sSpecial = @{366="TestEventUINT8Reliable"; 367="TestEventUINT8"; 368="TestEventUINT8Multicast"; 370="TestEventUINT8Multicast"}

# Rows where column T gets a specific test-function name instead of the
# default "x" marker.
$tSpecial = @{279="TestFieldUINT8"; 280="TestFieldUINT8"; 345="TestFieldUTF8DynamicReliable"}

foreach ($r in $rows) {
    $sCell = $ws.Cells.Item($r, 19)
    $tCell = $ws.Cells.Item($r, 20)

    if ($sSpecial.ContainsKey($r)) {
        $sCell.Value = $sSpecial[$r]
    } else {
        $sCell.Value = "x"
        $sCell.HorizontalAlignment = 1
    }

    if ($tSpecial.ContainsKey($r)) {
        $tCell.Value = $tSpecial[$r]
    } else {
        $tCell.Value = "x"
        $tCell.HorizontalAlignment = 1
    }
}

# Widen columns S and T to fit the newly added text (author applied a
# "best fit" auto-size after filling the columns in).
$ws.Columns.Item(19).ColumnWidth = 24.1
$ws.Columns.Item(20).ColumnWidth = 28.85

# Restore the view: scroll the frozen pane back to the top and select P12
# (matches the reviewer's cursor position after the edit).
$ws.Activate() | Out-Null
$ws.Range("P12").Select() | Out-Null

Write-Host "Applied S/T review columns to" $rows.Count "rows"
